# Re-generate the quadratic/linear experiment problem data.
# This mirrors the commit "volver a generar problemas cuadraticos y lineales":
# the leader/follower restriction expressions and the numeric vectors
# (Punto_modificado, Vector_bf, Vector_BF, Vector_Alpha) are refreshed with
# a freshly generated set of values, while the sheet layout itself is
# unchanged.

$wb = $excel.ActiveWorkbook

# --- Restricciones_del_follower ---------------------------------------
$ws3 = $wb.Worksheets.Item("Restricciones_del_follower")

$ws3.Range("A2").Value = "4.148236775818639 - x - 0.19647355163727953y_1 + 0.37783375314861467y_2"
$ws3.Range("B2").Value = "-4.148236775818639"
$ws3.Range("D2").Value = "0.34"
$ws3.Range("E2").Value = "5.1"
$ws3.Range("F2").Value = "0"

$ws3.Range("A3").Value = "-5.6381701630113366 + 0.9962378027778391y_1 + 1.0380436574444762y_2"
$ws3.Range("B3").Value = "5.6381701630113366"
$ws3.Range("D3").Value = "0.14"
$ws3.Range("E3").Value = "0"
$ws3.Range("F3").Value = "4.1"

$ws3.Range("A4").Value = "-2.295340050377834 + 0.4093198992443325y_1 - 0.7871536523929471y_2"
$ws3.Range("B4").Value = "-2.295340050377834"
$ws3.Range("D4").Value = "0.38"
$ws3.Range("E4").Value = "6.800000000000001"
$ws3.Range("F4").Value = "0"

# --- Punto_modificado ---------------------------------------------------
$ws4 = $wb.Worksheets.Item("Punto_modificado")
$ws4.Range("A2").Value = "5.25"
$ws4.Range("B2").Value = "1.7000000000000002"
$ws4.Range("C2").Value = "3.8"

# --- Vector_bf ------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Vector_bf")
$ws5.Range("A2").Value = "-1.928213846545069"
$ws5.Range("A3").Value = "-0.9746712002034358"

# --- Vector_BF ------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Vector_BF")
$ws6.Range("A2").Value = "5.1"
$ws6.Range("A3").Value = "-3.731360201511336"
$ws6.Range("A4").Value = "1.6256926952141062"

# --- Vector_Alpha -----------------------------------------------------------
$ws7 = $wb.Worksheets.Item("Vector_Alpha")
$ws7.Range("A2").Value = 2.25
$ws7.Range("A3").Value = 1.17
